$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 gets a new diary entry: date "01.11.18", begin 17:00, end 18:30,
# and a remark about the new Point-angle calculation feature.

# A20: enter the date as literal text (not an actual date value) so it
# keeps the worksheet's existing "date" cell style (numFmtId 164) instead
# of Excel re-formatting the cell as a real date/quoted-text style. We do
# this by writing it as a text formula first (which guarantees text type
# without touching the cell's number format), then flattening the formula
# down to its static value via copy / paste-special-values.
$ws.Range("A20").Formula = "=""01.11.18"""
$ws.Range("A20").Copy()
$ws.Range("A20").PasteSpecial(-4163)  # xlPasteValues

# B20 / C20: begin / end times (17:00 and 18:30 as Excel time fractions)
$ws.Range("B20").Value = 0.70833333333333337
$ws.Range("C20").Value = 0.77083333333333337

# E20: remark / activity description
$ws.Range("E20").Value = "-Add Point angle calculation"

# Reflect the new active cell/selection (matches where the user ended up)
$ws.Range("E20").Select()
